$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 126
$ws.Cells.Item(126, 2).Value = 0
$ws.Cells.Item(126, 3).Value = 32.07978307
$ws.Cells.Item(126, 4).Value = 0
$ws.Cells.Item(126, 6).Value = 6.80280734
$ws.Cells.Item(126, 7).Value = 5.04657086
$ws.Cells.Item(126, 8).Value = 3.2209418
$ws.Cells.Item(126, 9).Value = 13.94355935
$ws.Cells.Item(126, 10).Value = 55.36038006
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 6.96277437
$ws.Cells.Item(126, 13).Value = 7.85478858
$ws.Cells.Item(126, 15).Value = 23.01027401
$ws.Cells.Item(126, 16).Value = 12.73224361
$ws.Cells.Item(126, 17).Value = 3.91673055
$ws.Cells.Item(126, 18).Value = 12.50150664
$ws.Cells.Item(126, 19).Value = 2.63621836
$ws.Cells.Item(126, 20).Value = 14.52614555
$ws.Cells.Item(126, 21).Value = 12.64333675
$ws.Cells.Item(126, 22).Value = 17.19320619
$ws.Cells.Item(126, 23).Value = 8.553894509999999
$ws.Cells.Item(126, 24).Value = 13.15363959
$ws.Cells.Item(126, 25).Value = 0
$ws.Cells.Item(126, 26).Value = 11.56076124
$ws.Cells.Item(126, 27).Value = 7.41078836
$ws.Cells.Item(126, 28).Value = 14.774324
$ws.Cells.Item(126, 30).Value = 5.18150865
$ws.Cells.Item(126, 31).Value = 0
$ws.Cells.Item(126, 32).Value = 11.17758416
$ws.Cells.Item(126, 33).Value = 91.83488593
$ws.Cells.Item(126, 34).Value = 23.55689551
$ws.Cells.Item(126, 35).Value = 0
$ws.Cells.Item(126, 36).Value = 4.6159519
$ws.Cells.Item(126, 37).Value = 13.54985888
$ws.Cells.Item(126, 38).Value = 12.38255236
$ws.Cells.Item(126, 39).Value = 8.052551920000001
$ws.Cells.Item(126, 40).Value = 7.37004402
$ws.Cells.Item(126, 41).Value = 0.71876079
$ws.Cells.Item(126, 42).Value = 0
$ws.Cells.Item(126, 43).Value = 11.46008882
$ws.Cells.Item(126, 45).Value = 3.51636536
$ws.Cells.Item(126, 46).Value = 4.7534468
$ws.Cells.Item(126, 47).Value = 16.5254644
$ws.Cells.Item(126, 48).Value = 11.89476536
$ws.Cells.Item(126, 49).Value = 6.55495944
$ws.Cells.Item(126, 50).Value = 30.01355252
$ws.Cells.Item(126, 51).Value = 8.557738369999999
$ws.Cells.Item(126, 53).Value = 51.55977704
$ws.Cells.Item(126, 54).Value = 11.3179946
$ws.Cells.Item(126, 55).Value = 8.955335959999999
$ws.Cells.Item(126, 56).Value = 11.48459981
$ws.Cells.Item(126, 57).Value = 0

# Row 127
$ws.Cells.Item(127, 2).Value = 12.11784292
$ws.Cells.Item(127, 3).Value = 12.90241976
$ws.Cells.Item(127, 4).Value = 11.06935687
$ws.Cells.Item(127, 6).Value = 19.02024904
$ws.Cells.Item(127, 7).Value = 4.64298119
$ws.Cells.Item(127, 8).Value = 13.62890831
$ws.Cells.Item(127, 9).Value = 15.43696148
$ws.Cells.Item(127, 10).Value = 17.58095808
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 7.48620792
$ws.Cells.Item(127, 13).Value = 5.20524368
$ws.Cells.Item(127, 15).Value = 24.6347654
$ws.Cells.Item(127, 16).Value = 33.14337674
$ws.Cells.Item(127, 17).Value = 0
$ws.Cells.Item(127, 18).Value = 10.81020647
$ws.Cells.Item(127, 19).Value = 8.342476420000001
$ws.Cells.Item(127, 20).Value = 7.67983638
$ws.Cells.Item(127, 21).Value = 7.05930321
$ws.Cells.Item(127, 22).Value = 12.1142975
$ws.Cells.Item(127, 23).Value = 4.7618873
$ws.Cells.Item(127, 24).Value = 11.31835808
$ws.Cells.Item(127, 25).Value = 0
$ws.Cells.Item(127, 26).Value = 4.67643394
$ws.Cells.Item(127, 27).Value = 17.36561475
$ws.Cells.Item(127, 28).Value = 5.73872636
$ws.Cells.Item(127, 30).Value = 16.54105508
$ws.Cells.Item(127, 31).Value = 0
$ws.Cells.Item(127, 32).Value = 6.13644652
$ws.Cells.Item(127, 33).Value = 0
$ws.Cells.Item(127, 34).Value = 21.61167123
$ws.Cells.Item(127, 35).Value = 23.81508638
$ws.Cells.Item(127, 36).Value = 11.44296435
$ws.Cells.Item(127, 37).Value = 13.97531241
$ws.Cells.Item(127, 38).Value = 4.83495465
$ws.Cells.Item(127, 39).Value = 6.97583575
$ws.Cells.Item(127, 40).Value = 7.97452242
$ws.Cells.Item(127, 41).Value = 14.12392325
$ws.Cells.Item(127, 42).Value = 0.58220671
$ws.Cells.Item(127, 43).Value = 11.1920124
$ws.Cells.Item(127, 45).Value = 3.55089838
$ws.Cells.Item(127, 46).Value = 7.52339399
$ws.Cells.Item(127, 47).Value = 33.73535995
$ws.Cells.Item(127, 48).Value = 5.73327049
$ws.Cells.Item(127, 49).Value = 5.86441641
$ws.Cells.Item(127, 50).Value = 23.22632091
$ws.Cells.Item(127, 51).Value = 11.31163493
$ws.Cells.Item(127, 53).Value = 0
$ws.Cells.Item(127, 54).Value = 2.0799488
$ws.Cells.Item(127, 55).Value = 6.43477271
$ws.Cells.Item(127, 56).Value = 4.14052513
$ws.Cells.Item(127, 57).Value = 7.78016808

# Row 128
$ws.Cells.Item(128, 2).Value = 0
$ws.Cells.Item(128, 3).Value = 20.1847512
$ws.Cells.Item(128, 4).Value = 17.36458959
$ws.Cells.Item(128, 6).Value = 11.97709067
$ws.Cells.Item(128, 7).Value = 4.23053211
$ws.Cells.Item(128, 8).Value = 4.5700691
$ws.Cells.Item(128, 9).Value = 13.62297787
$ws.Cells.Item(128, 10).Value = 33.41345571
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 4.21639424
$ws.Cells.Item(128, 13).Value = 6.12105556
$ws.Cells.Item(128, 15).Value = 14.40135234
$ws.Cells.Item(128, 16).Value = 9.12511645
$ws.Cells.Item(128, 17).Value = 0
$ws.Cells.Item(128, 18).Value = 9.35806429
$ws.Cells.Item(128, 19).Value = 1.30880036
$ws.Cells.Item(128, 20).Value = 5.29410722
$ws.Cells.Item(128, 21).Value = 13.623206
$ws.Cells.Item(128, 22).Value = 13.79106166
$ws.Cells.Item(128, 23).Value = 6.4331438
$ws.Cells.Item(128, 24).Value = 9.23574571
$ws.Cells.Item(128, 25).Value = 0
$ws.Cells.Item(128, 26).Value = 10.45310043
$ws.Cells.Item(128, 27).Value = 8.715119769999999
$ws.Cells.Item(128, 28).Value = 8.527367979999999
$ws.Cells.Item(128, 30).Value = 9.297101469999999
$ws.Cells.Item(128, 31).Value = 0
$ws.Cells.Item(128, 32).Value = 6.56132347
$ws.Cells.Item(128, 33).Value = 58.9119578
$ws.Cells.Item(128, 34).Value = 13.13776974
$ws.Cells.Item(128, 35).Value = 22.69229168
$ws.Cells.Item(128, 36).Value = 6.30727045
$ws.Cells.Item(128, 37).Value = 9.18140395
$ws.Cells.Item(128, 38).Value = 5.2439938
$ws.Cells.Item(128, 39).Value = 6.45307393
$ws.Cells.Item(128, 40).Value = 5.63611264
$ws.Cells.Item(128, 41).Value = 1.97619827
$ws.Cells.Item(128, 42).Value = 0
$ws.Cells.Item(128, 43).Value = 9.076702640000001
$ws.Cells.Item(128, 45).Value = 7.25985526
$ws.Cells.Item(128, 46).Value = 6.63328624
$ws.Cells.Item(128, 47).Value = 6.34905301
$ws.Cells.Item(128, 48).Value = 6.40174806
$ws.Cells.Item(128, 49).Value = 5.11233129
$ws.Cells.Item(128, 50).Value = 17.08847955
$ws.Cells.Item(128, 51).Value = 7.02619867
$ws.Cells.Item(128, 53).Value = 29.86337551
$ws.Cells.Item(128, 54).Value = 8.506800500000001
$ws.Cells.Item(128, 55).Value = 5.78040292
$ws.Cells.Item(128, 56).Value = 20.96817893
$ws.Cells.Item(128, 57).Value = 0

# New date rows
$ws.Cells.Item(131, 1).Value = "09 06 2020"
$ws.Cells.Item(132, 1).Value = "10 06 2020"
$ws.Cells.Item(133, 1).Value = "11 06 2020"
